# Update countries & provincias Spain
# Applies the 22-Jul-2020 16:50 data refresh to the "Pais" sheet:
#  - refreshed case/death/recovered counters for a number of countries
#  - three countries swapped rank (and therefore row) with their neighbour:
#      Tunez / Namibia        (rows 136-137)
#      Mauricio / Birmania    (rows 164-165)
#      Groenlandia / Islas Malvinas (rows 210-211)
#  - the "datos actualizados" timestamp label was bumped

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- timestamp label -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 22 de Julio de 2020 a las 16:50"

# --- plain numeric refreshes ------------------------------------------
$ws.Range("B4").Value = 4030936
$ws.Range("C4").Value = 2367
$ws.Range("D4").Value = 1888420
$ws.Range("E4").Value = 1997529
$ws.Range("G4").Value = 34
$ws.Range("H4").Value = 144987

$ws.Range("B5").Value = 2167988
$ws.Range("C5").Value = 1456
$ws.Range("E5").Value = 620390
$ws.Range("G5").Value = 31
$ws.Range("H5").Value = 81628

$ws.Range("B6").Value = 1216965
$ws.Range("C6").Value = 22880
$ws.Range("D6").Value = 769979
$ws.Range("E6").Value = 417512
$ws.Range("G6").Value = 704
$ws.Range("H6").Value = 29474

$ws.Range("B21").Value = 204153
$ws.Range("C21").Value = 263
$ws.Range("E21").Value = 6873

$ws.Range("D46").Value = 44795
$ws.Range("E46").Value = 3922

$ws.Range("B63").Value = 21798
$ws.Range("C63").Value = 356
$ws.Range("E63").Value = 6230
$ws.Range("G63").Value = 5
$ws.Range("H63").Value = 712

$ws.Range("B91").Value = 7015
$ws.Range("C91").Value = 48
$ws.Range("D91").Value = 5741
$ws.Range("E91").Value = 1216
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 58

$ws.Range("B114").Value = 2731
$ws.Range("C114").Value = 1
$ws.Range("E114").Value = 656

$ws.Range("B115").Value = 2494
$ws.Range("C115").Value = 17
$ws.Range("D115").Value = 1869
$ws.Range("E115").Value = 502
$ws.Range("G115").Value = 1
$ws.Range("H115").Value = 123

$ws.Range("B116").Value = 2462
$ws.Range("C116").Value = 13
$ws.Range("D116").Value = 2326
$ws.Range("E116").Value = 49

$ws.Range("B158").Value = 561
$ws.Range("C158").Value = 21
$ws.Range("D158").Value = 165
$ws.Range("E158").Value = 364
$ws.Range("G158").Value = 1
$ws.Range("H158").Value = 32

# --- Tunez / Namibia swap (rows 136-137) ------------------------------
$ws.Range("A136").Value = "Namibia"
$ws.Range("B136").Value = 1402
$ws.Range("C136").Value = 36
$ws.Range("D136").Value = 64
$ws.Range("E136").Value = 1331
$ws.Range("H136").Value = 7

$ws.Range("A137").Value = "Tunez"
$ws.Range("B137").Value = 1389
$ws.Range("D137").Value = 1103
$ws.Range("E137").Value = 236
$ws.Range("H137").Value = 50

# --- Mauricio / Birmania swap (rows 164-165) --------------------------
$ws.Range("A164").Value = "Birmania"
$ws.Range("C164").Value = 2
$ws.Range("D164").Value = 280
$ws.Range("E164").Value = 57
$ws.Range("H164").Value = 6

$ws.Range("A165").Value = "Mauricio"
$ws.Range("B165").Value = 343
$ws.Range("D165").Value = 332
$ws.Range("E165").Value = 1
$ws.Range("H165").Value = 10

# --- Groenlandia / Islas Malvinas swap (rows 210-211) ------------------
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("A211").Value = "Groenlandia"
